$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the orphan "5840560 - Marco Antonio Carvalho Pereira" data row under
# "Docentes responsaveis:") - this shifts rows 14-24 up by one.
$ws.Rows.Item(13).Delete()

# Replace the "Objetivos:" answer (row 10) with the docente text.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Replace the "Programa resumido:" answer (row 13 after shift) with "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Replace the "Programa:" answer (row 15 after shift) with the activation date.
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

# Replace the "Método:" answer (row 18 after shift) with the docente text.
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Replace the "Critério:" answer (row 19 after shift) with the teaching method text.
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas de exercícios."

# Replace the "Norma de recuperação:" answer (row 20 after shift) with the evaluation text.
$ws.Range("B20").Value = "Média de Provas e trabalhos (MF)."
$ws.Range("C20").Value = "Média de Provas e trabalhos (MF)."

# Replace the "Bibliografia:" answer (row 21 after shift) with the recovery-exam text.
$ws.Range("B21").Value = "Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR"
$ws.Range("C21").Value = "Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR"
